# Update the "Förändrad" (Changed) date column (C) for rows 2-24
# from serial date 45226 (2023-10-27) to 45227 (2023-10-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 24; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45226) {
        $cell.Value = 45227
    }
}
